$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("E2").Value = 29740.57209542279
$ws.Range("G2").Value = 8095.925712661859
$ws.Range("I2").Value = 14200.72657838
$ws.Range("L2").Value = 53380.29012525
$ws.Range("M2").Value = 10658.387169815
$ws.Range("N2").Value = 7637.402481230281
$ws.Range("O2").Value = 7117.919615036673

# Sheet "2030" (sheet2)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 3208.898512992542
$ws.Range("E2").Value = 47193.8449614085
$ws.Range("G2").Value = 8095.925712661859
$ws.Range("I2").Value = 23823.32840570365
$ws.Range("L2").Value = 75872.38143273753
$ws.Range("M2").Value = 16509.210294786
$ws.Range("N2").Value = 9395.063137084702
$ws.Range("O2").Value = 8358.401204250757

# Sheet "2035" (sheet3)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 1553.56866473945
$ws.Range("B2").Value = 5378.646016954563
$ws.Range("E2").Value = 60525.83365072873
$ws.Range("G2").Value = 8095.925712661859
$ws.Range("I2").Value = 41536.04551959009
$ws.Range("L2").Value = 75872.38143273753
$ws.Range("M2").Value = 22090.099315812
$ws.Range("N2").Value = 13763.87675960632
$ws.Range("O2").Value = 12838.44160168211

# Sheet "2040" (sheet4)
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 1553.56866473945
$ws.Range("B2").Value = 5378.646016954563
$ws.Range("E2").Value = 60525.83365072873
$ws.Range("G2").Value = 8095.925712661859
$ws.Range("I2").Value = 41536.04551959009
$ws.Range("L2").Value = 75872.38143273753
$ws.Range("M2").Value = 22090.099315812
$ws.Range("N2").Value = 13763.87675960632
$ws.Range("O2").Value = 12838.44160168211

# Sheet "2045" (sheet5)
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 1553.56866473945
$ws.Range("B2").Value = 5378.646016954563
$ws.Range("E2").Value = 60525.83365072873
$ws.Range("G2").Value = 8095.925712661859
$ws.Range("I2").Value = 41536.04551959009
$ws.Range("L2").Value = 75872.38143273753
$ws.Range("M2").Value = 22090.099315812
$ws.Range("N2").Value = 13763.87675960632
$ws.Range("O2").Value = 12838.44160168211

# Sheet "2050" (sheet6)
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 1553.56866473945
$ws.Range("B2").Value = 5378.646016954563
$ws.Range("E2").Value = 60525.83365072873
$ws.Range("G2").Value = 8095.925712661859
$ws.Range("I2").Value = 41536.04551959009
$ws.Range("L2").Value = 75872.38143273753
$ws.Range("M2").Value = 22090.099315812
$ws.Range("N2").Value = 13763.87675960632
$ws.Range("O2").Value = 12838.44160168211
